$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 1176-1177, shifting existing rows 1176-1231 down to 1178-1233
$ws.Range("A1176:A1177").EntireRow.Insert()

# Row 1176
$ws.Range("A1176").Value = 6
$ws.Range("B1176").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1176").Value = "Metropolitana"
$ws.Range("D1176").Value = 45041
$ws.Range("E1176").Value = 13
$ws.Range("F1176").Value = 100112023
$ws.Range("G1176").Value = "Brócoli"
$ws.Range("H1176").Value = "Sin especificar"
$ws.Range("I1176").Value = "Primera"
$ws.Range("J1176").Value = 11800
$ws.Range("K1176").Value = 800
$ws.Range("L1176").Value = 900
$ws.Range("M1176").Value = 845
$ws.Range("N1176").Value = "`$/unidad"
$ws.Range("O1176").Value = "Región Metropolitana"
$ws.Range("P1176").Value = 845
$ws.Range("Q1176").Value = 1
$ws.Range("R1176").Value = "Hortaliza"

# Row 1177
$ws.Range("A1177").Value = 6
$ws.Range("B1177").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1177").Value = "Metropolitana"
$ws.Range("D1177").Value = 45041
$ws.Range("E1177").Value = 13
$ws.Range("F1177").Value = 100112023
$ws.Range("G1177").Value = "Brócoli"
$ws.Range("H1177").Value = "Sin especificar"
$ws.Range("I1177").Value = "Segunda"
$ws.Range("J1177").Value = 4600
$ws.Range("K1177").Value = 500
$ws.Range("L1177").Value = 600
$ws.Range("M1177").Value = 578
$ws.Range("N1177").Value = "`$/unidad"
$ws.Range("O1177").Value = "Región Metropolitana"
$ws.Range("P1177").Value = 578
$ws.Range("Q1177").Value = 1
$ws.Range("R1177").Value = "Hortaliza"

# Row 1232
$ws.Range("A1232").Value = 6
$ws.Range("B1232").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1232").Value = "Metropolitana"
$ws.Range("D1232").Value = 45007
$ws.Range("E1232").Value = 13
$ws.Range("F1232").Value = 100112023
$ws.Range("G1232").Value = "Brócoli"
$ws.Range("H1232").Value = "Sin especificar"
$ws.Range("I1232").Value = "Primera"
$ws.Range("J1232").Value = 10200
$ws.Range("K1232").Value = 700
$ws.Range("L1232").Value = 800
$ws.Range("M1232").Value = 746
$ws.Range("N1232").Value = "`$/unidad"
$ws.Range("O1232").Value = "Región Metropolitana"
$ws.Range("P1232").Value = 746
$ws.Range("Q1232").Value = 1
$ws.Range("R1232").Value = "Hortaliza"

# Row 1233
$ws.Range("A1233").Value = 6
$ws.Range("B1233").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C1233").Value = "Metropolitana"
$ws.Range("D1233").Value = 45007
$ws.Range("E1233").Value = 13
$ws.Range("F1233").Value = 100112023
$ws.Range("G1233").Value = "Brócoli"
$ws.Range("H1233").Value = "Sin especificar"
$ws.Range("I1233").Value = "Segunda"
$ws.Range("J1233").Value = 2600
$ws.Range("K1233").Value = 600
$ws.Range("L1233").Value = 600
$ws.Range("M1233").Value = 600
$ws.Range("N1233").Value = "`$/unidad"
$ws.Range("O1233").Value = "Región Metropolitana"
$ws.Range("P1233").Value = 600
$ws.Range("Q1233").Value = 1
$ws.Range("R1233").Value = "Hortaliza"
